$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest values.
# NumberFormat is set to Text ("@") first so the numeric-looking / percent-looking
# strings are preserved as literal text, matching the source data (t="inlineStr").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.36%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.61%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.712"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.57%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.76%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.836"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.13%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.509"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.01%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.987"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.97%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.13%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9278"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.80%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1254"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.16%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1962"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.94%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09423"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.15%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03996"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "9.25%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.84%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001299"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.30%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006096"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.98%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.438"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.64%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.124"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "9.61%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1373"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.00%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2627"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.78%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04424"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.03%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001242"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.38%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.01%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.88%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003995"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.06%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02815"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.08%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05537"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.64%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007924"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.88%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1439"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.04%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008964"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.82%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.60%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01038"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.16%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007142"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.86%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.21%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003449"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "18.06%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002281"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.18%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.21%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.21%"
